$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record (Fecha 44783 = 2022-08-10) is inserted right above the
# existing row 257 ("Start Ruby" / "Primera", 110 @ 12000, $/bandeja 15 kilos
# granel, 800 $/Kg, 15 Kg/unidad). Inserting a whole row shifts every
# following row (old 257..276) down by one (to 258..277), matching the diff.
$ws.Rows.Item(257).Insert()

# Populate the newly-inserted row 257 with the same data as row 256 (which
# sits right above it and is otherwise untouched by the edit), then correct
# its date to the new value.
$ws.Range("A256:T256").Copy()
$ws.Range("A257").PasteSpecial()
$ws.Range("D257").Value = 44783
